# Apply odds updates to Sheet1 as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 2.5
$ws.Range("I2").Value = 3.25
$ws.Range("J2").Value = 3.5
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 6
$ws.Range("U2").Value = 2.38
$ws.Range("V2").Value = 1.53
$ws.Range("W2").Value = 5.5
$ws.Range("AD2").Value = 6
$ws.Range("AI2").Value = 13
$ws.Range("AU2").Value = 10

# Row 3
$ws.Range("G3").Value = 1.57
$ws.Range("H3").Value = 4.2
$ws.Range("I3").Value = 5.25
$ws.Range("S3").Value = 1.33
$ws.Range("T3").Value = 3.25
$ws.Range("W3").Value = 8
$ws.Range("X3").Value = 8
$ws.Range("AB3").Value = 23
$ws.Range("AD3").Value = 8.5
$ws.Range("AE3").Value = 17
$ws.Range("AF3").Value = 51
$ws.Range("AH3").Value = 15
$ws.Range("AN3").Value = 3.6
$ws.Range("AT3").Value = 3.25

# Row 4
$ws.Range("G4").Value = 1.9
$ws.Range("H4").Value = 3.3
$ws.Range("I4").Value = 4.33
$ws.Range("J4").Value = 2.6
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 8.5
$ws.Range("Z4").Value = 15
$ws.Range("AO4").Value = 10
$ws.Range("AQ4").Value = 34

$wb.Save()
